$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes
# Note: the COM layer quantizes ColumnWidth to the nearest 1/6 (pixel) boundary,
# so the inputs below are chosen to reliably land on the closest achievable
# raw width to the target values from the diff.
$ws.Columns.Item(3).ColumnWidth = 2.333
$ws.Columns.Item(5).ColumnWidth = 2.333
$ws.Columns.Item(11).ColumnWidth = 1.333
$ws.Columns.Item(17).ColumnWidth = 4.833

# Row 1 value changes
$ws.Range("C1").Value = 25
$ws.Range("E1").Value = 19
$ws.Range("F1").Value = 23
$ws.Range("G1").Value = 31
$ws.Range("I1").Value = 26
$ws.Range("J1").Value = 31
$ws.Range("K1").Value = 5
$ws.Range("L1").Value = 33
$ws.Range("M1").Value = 0.068999999999999992
$ws.Range("N1").Value = 0.087999999999999995
$ws.Range("O1").Value = 0.060999999999999999
$ws.Range("P1").Value = 0.074999999999999997
$ws.Range("Q1").Value = 0.086999999999999994
